# Update scripts with new TPM values.
# The rows whose Sending cluster was "ECs" (original rows 2-4) are removed,
# which shifts the former "FAPs" (rows 5-7) and "MuSCs" (rows 8-10) sender
# rows up to become the new rows 2-4 and 5-7. The remaining rows keep their
# original text labels but get refreshed numeric (TPM-derived) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows belonging to the "ECs" sending cluster.
$ws.Rows("2:4").Delete()

# New numeric data (columns E through T) for the remaining rows, now
# located at rows 2-7 after the deletion above.
$data = @{
    2 = @(2, 0.6666666666666666, 1.154926666666667, 3.46478, 0.9439377097772548, 0.9439377097772548, 3, 1, 2.027115333333333, 6.081346, 0.006596284565418616, 0.006596284565418615, 2.341169554875555, 21.07052599388, 0.006226481745720302, 0.006226481745720301)
    3 = @(2, 0.6666666666666666, 1.154926666666667, 3.46478, 0.9439377097772548, 0.9439377097772548, 3, 1, 256.4443053333333, 769.332916, 0.8344762556643375, 0.8344762556643374, 296.1743667442755, 2665.56930069848, 0.7876936056352937, 0.7876936056352936)
    4 = @(2, 0.6666666666666666, 1.154926666666667, 3.46478, 0.9439377097772548, 0.9439377097772548, 3, 1, 48.84026566666667, 146.520797, 0.158927459770244, 0.158927459770244, 56.40692522551779, 507.6623270296601, 0.1500176223962409, 0.1500176223962409)
    5 = @(2, 0.6666666666666666, 0.06859333333333334, 0.20578, 0.05606229022274531, 0.05606229022274531, 3, 1, 2.027115333333333, 6.081346, 0.006596284565418616, 0.006596284565418615, 0.1390465977644444, 1.25141937988, 0.0003698028196983139, 0.0003698028196983138)
    6 = @(2, 0.6666666666666666, 0.06859333333333334, 0.20578, 0.05606229022274531, 0.05606229022274531, 3, 1, 256.4443053333333, 769.332916, 0.8344762556643375, 0.8344762556643374, 17.59036971716445, 158.31332745448, 0.0467826500290439, 0.04678265002904389)
    7 = @(2, 0.6666666666666666, 0.06859333333333334, 0.20578, 0.05606229022274531, 0.05606229022274531, 3, 1, 48.84026566666667, 146.520797, 0.158927459770244, 0.158927459770244, 3.350116622962223, 30.15104960666001, 0.0089098373740031, 0.008909837374003098)
}

$columns = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $ws.Range("$col$rowNum").Value = $values[$i]
    }
}
